# Edit script for "(宣道詩279) 耶穌我救主降生伯利恆.pptx"
#
# Summary of the change being applied:
#   1. Slide 1 (title) and slide 2 (content placeholder) each had the
#      song title split across two runs ("耶" + "穌我救主降生伯利恆");
#      they are merged back into a single run.
#   2. Every "page counter" textbox ("( 1 )", "( 2 )", "( 3 )", "( 4 )")
#      is resized (shorter height), its font shrunk from 48pt to 32pt,
#      and its text changed to show "( n / 4 )" instead of just "( n )".
#      On the first occurrence of each number (slides 2, 6, 10, 14) the
#      text is produced as three separate runs (mirroring how PowerPoint
#      split the run when the new text was typed in the middle of it).
#      On slide 2 specifically the runs also pick up an explicit
#      "+mn-lt" (theme minor-latin) Latin typeface.

$p = $ppt.ActivePresentation

$NEW_HEIGHT_PT = 584775 / 12700   # EMU -> points
$NEW_SIZE = 32

function Merge-TitleRuns($shape) {
    # The title text "耶" + "穌我救主降生伯利恆" live as two consecutive
    # runs at the very end of the text range; re-writing that span as a
    # single string collapses them into one run (keeping the formatting
    # of the first of the two runs), exactly like the authored edit.
    $tr = $shape.TextFrame.TextRange
    $whole = $tr.Text
    $idx = $whole.IndexOf("耶穌我救主降生伯利恆")
    $seg = $tr.Characters($idx + 1, 10)
    $seg.Text = "耶穌我救主降生伯利恆"
}

function Update-PageCounter($shape, [string]$num, [bool]$split3, [bool]$mnlt) {
    # Shrink the textbox height.
    $shape.Height = $NEW_HEIGHT_PT

    $tr = $shape.TextFrame.TextRange
    $tr.Font.Size = $NEW_SIZE
    if ($mnlt) {
        $tr.Font.Name = "+mn-lt"
    }

    if ($split3 -and $mnlt) {
        # slide 2 special case: "( N )" -> "( " / "N / 4 " / ")"  (three runs)
        $mid = $tr.Characters(3, 1)
        $mid.Text = "$num / 4 "
        $tail = $tr.Characters($tr.Length - 1, 2)
        $tail.Text = ")"
    } elseif ($split3) {
        # "( N )" -> "( " / "N " / "/ 4 )"  (three runs)
        $mid = $tr.Characters(3, 1)
        $mid.Text = "$num "
        $tail = $tr.Characters(5, 2)
        $tail.Text = "/ 4 )"
    } else {
        # "( N )" -> "( N / 4 )" (single run)
        $tr.Text = "( $num / 4 )"
    }
}

# ---- Slide 1 : title only ------------------------------------------------
Merge-TitleRuns $p.Slides.Item(1).Shapes.Item(1)

# ---- Slide 2 : title merge + "(1)" -> "(1 / 4)" with split runs + +mn-lt -
$s2 = $p.Slides.Item(2)
Merge-TitleRuns $s2.Shapes.Item(1)
Update-PageCounter $s2.Shapes.Item(2) "1" $true $true

# ---- Slides 3-5 : "(1)" -> "(1 / 4)" -------------------------------------
foreach ($i in 3,4,5) {
    Update-PageCounter $p.Slides.Item($i).Shapes.Item(2) "1" $false $false
}

# ---- Slide 6 : "(2)" -> "(2 / 4)" with split runs ------------------------
Update-PageCounter $p.Slides.Item(6).Shapes.Item(2) "2" $true $false

# ---- Slides 7-9 : "(2)" -> "(2 / 4)" --------------------------------------
foreach ($i in 7,8,9) {
    Update-PageCounter $p.Slides.Item($i).Shapes.Item(2) "2" $false $false
}

# ---- Slide 10 : "(3)" -> "(3 / 4)" with split runs -----------------------
Update-PageCounter $p.Slides.Item(10).Shapes.Item(2) "3" $true $false

# ---- Slides 11-13 : "(3)" -> "(3 / 4)" -------------------------------------
foreach ($i in 11,12,13) {
    Update-PageCounter $p.Slides.Item($i).Shapes.Item(2) "3" $false $false
}

# ---- Slide 14 : "(4)" -> "(4 / 4)" with split runs -----------------------
Update-PageCounter $p.Slides.Item(14).Shapes.Item(2) "4" $true $false

# ---- Slides 15-17 : "(4)" -> "(4 / 4)" -------------------------------------
foreach ($i in 15,16,17) {
    Update-PageCounter $p.Slides.Item($i).Shapes.Item(2) "4" $false $false
}
